# Add floating rate handling: append new interest-rate rows (2026-12-31 curve)
# to the "interest rates" sheet, mirroring the existing 2022-12-31 block with
# rates shifted up by 2% (0.02).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("interest rates")

# New rows 12-20 mirror rows 3-11 (same Name/Type/Tenor/Maturity labels),
# but with a new date and rates increased by 0.02.
$newDate = 46387  # 2026-12-31 (serial date, same format as existing rows)

$rows = @(
    @{ Row=12; B="Euribor"; C="Spot"; D="3m";  E=$null; F=0.0505 },
    @{ Row=13; B="Euribor"; C="Spot"; D="6m";  E=$null; F=0.0495 },
    @{ Row=14; B="Euribor"; C="Zero"; D=$null; E="1m";  F=0.0518 },
    @{ Row=15; B="Euribor"; C="Zero"; D=$null; E="1y";  F=0.048600000000000004 },
    @{ Row=16; B="Euribor"; C="Zero"; D=$null; E="10y"; F=0.0455 },
    @{ Row=17; B="Euribor"; C="Zero"; D=$null; E="20y"; F=0.0465 },
    @{ Row=18; B="Euribor"; C="Zero"; D=$null; E="30y"; F=0.047 },
    @{ Row=19; B="Euribor"; C="Swap"; D="3M";  E="1Y";  F=0.0485 },
    @{ Row=20; B="Euribor"; C="Swap"; D="3M";  E="10Y"; F=0.0455 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column A: date, formatted/styled like the existing date column (copy style from A3)
    $ws.Range("A3").Copy($ws.Range("A" + $rowNum))
    $ws.Range("A" + $rowNum).Value = $newDate

    # Column B: Name (text)
    $ws.Range("B" + $rowNum).Value = $r.B

    # Column C: Type (text)
    $ws.Range("C" + $rowNum).Value = $r.C

    # Column D: Tenor (text, optional)
    if ($r.D -ne $null) {
        $ws.Range("D" + $rowNum).Value = $r.D
    }

    # Column E: Maturity (text, optional)
    if ($r.E -ne $null) {
        $ws.Range("E" + $rowNum).Value = $r.E
    }

    # Column F: Rate, styled like the existing rate column (copy style from F3)
    $ws.Range("F3").Copy($ws.Range("F" + $rowNum))
    $ws.Range("F" + $rowNum).Value = $r.F
}

# Update the active selection on the sheet to match the edited state.
$ws.Range("J13").Select() | Out-Null
